$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.741.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.315.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.49%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'246.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -6.06%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'648.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -14.90%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -11.35%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.16%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.964"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -14.51%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.316.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -7.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'39.47"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -8.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'96.578.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.77%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -5.74%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000248"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -9.72%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.934.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.51%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'8.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.64%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.313.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.56%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -6.32%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -8.26%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'BitcoinCash"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'493.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -6.98%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'Uniswap"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'10.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -6.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.445"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -9.65%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -10.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'6.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'93.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -9.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'11.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -9.33%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.495.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.28%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.139"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -9.18%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'10.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -9.57%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.184"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -6.88%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +6.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.53%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.536"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -8.73%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'27.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -9.65%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +1.68%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -8.26%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.149"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -8.43%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'498.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -7.89%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'24.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.14%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.33%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.814"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -6.90%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0400"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -11.19%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'8.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.18%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'5.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.74%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'52.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.48%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -11.99%  "
$ws.Range("E51").Style = "Normal"
